# Release R01.00.58-20150625: log a new entry for TestShell/Modules/audit.js
# (timer task of audit not implemented on server side, so the related code
# was commented out).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Append the new log row right after the existing last row (46).
$ws.Range("B47").Value = "TestShell/Modules/audit.js"
$ws.Range("C47").Value = "Commented out line 106-116, timer task of audit is not implemented on server side"

# Scroll the view down to the new row and leave the selection just below it,
# matching where the author's cursor ended up after typing the new entry.
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("C48").Select()
